$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
SELECT
    COUNT(DISTINCT std.study_ID) AS "Studies",
    COUNT(DISTINCT prt.participant_id) AS "Participants",
    COUNT(DISTINCT smp.sample_id) AS "Samples",
    (COUNT(DISTINCT seq.id) + COUNT(DISTINCT paf.id)) AS "Files"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_sequencing_file seq ON smp.id = seq."sample.id"
LEFT JOIN 
    df_pathology_file paf ON smp.id = paf."sample.id"
WHERE 
    std.study_ID = 'phs002430' 
    AND prt.race = 'Asian' 
    AND prt.sex_at_birth = 'Female';
'@

$ws.Range("C2").Value = $newQuery
$ws.Range("C2").Select()
